$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record is published. It becomes the new first data
# row (row 2), and every existing data row from 2-10 shifts down to 3-11,
# aging the previous row 11 record out of this 10-row rolling window
# (rows 12/13 below are a separate, untouched block of older history).
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

$ws.Range("A2").Value = 6
$ws.Range("B2").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C2").Value = "Metropolitana"
$ws.Range("D2").Value = 44453
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E2").Value = 13
$ws.Range("F2").Value = 100112035
$ws.Range("G2").Value = "Bruselas (repollito)"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 280
$ws.Range("K2").Value = 20000
$ws.Range("L2").Value = 22000
$ws.Range("M2").Value = 21286
$ws.Range("N2").Value = "$/malla 15 kilos"
$ws.Range("O2").Value = "Provincia de Quillota"
$ws.Range("P2").Value = 1419
$ws.Range("Q2").Value = 15
$ws.Range("R2").Value = "Hortaliza"

# The old row 11 (now pushed down to row 12 by the insert) drops out of
# the rolling window; remove it so the trailing historical rows keep
# their original row numbers / values.
$ws.Rows.Item(12).Delete()
